$d = $word.ActiveDocument

# Locate the run of text "En esta sección " (it is currently a single run
# with a trailing space) without disturbing anything else in the document.
$findRange = $d.Content.Duplicate
$found = $findRange.Find.Execute("En esta sección ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $runStart = $findRange.Start
    $runEnd = $findRange.End

    # That run is immediately followed by another run ("se ") that happens
    # to share identical run formatting (rFonts cs=Arial, noProof,
    # lang=es-PE). This engine automatically re-coalesces any text edited
    # inside "En esta sección " with that following run, which would wipe
    # out its distinguishing w:rsidR attribute even though the diff we
    # need to reproduce leaves that run completely untouched.
    #
    # To prevent that unwanted merge, we temporarily give the run a
    # distinguishing format (Bold) before editing its text; this stops the
    # coalescing pass from pulling the following "se " run into our edit.
    # We restore the formatting afterwards, and only then split the newly
    # written text into separate runs ("En esta sección" / "," / " ") by
    # toggling Bold on/off for each piece - without ever touching the
    # "se " run again.

    $run1 = $d.Range($runStart, $runEnd)
    $run1.Bold = 1

    # Insert ", " right before the following run (i.e. right after the
    # trailing space that currently ends "En esta sección ").
    $ins = $d.Range($runEnd, $runEnd)
    $ins.InsertBefore(", ")

    # Remove the original trailing space that used to follow "sección"
    # (it sits right before the text we just inserted).
    $oldSpace = $d.Range($runEnd - 1, $runEnd)
    $oldSpace.Text = ""

    # Restore normal (non-bold) formatting across the whole edited span:
    # "En esta sección, " (comma + space included, "se " excluded).
    $edited = $d.Range($runStart, $runEnd + 1)
    $edited.Bold = 0

    # Split the trailing "," off into its own run.
    $commaR = $d.Range($runEnd - 1, $runEnd)
    $commaR.Bold = 1
    $commaR.Bold = 0

    # Split the trailing " " off into its own run.
    $spaceR = $d.Range($runEnd, $runEnd + 1)
    $spaceR.Bold = 1
    $spaceR.Bold = 0
}
